# Generate Report for Handback
# Applies the "handback" updates to the localization-status workbook:
#   - Overview sheet status column wording changes from "Ready for handoff"
#     to "Handed back: in sync with en-US"
#   - Each per-locale sheet (zh-cn, de-de) gets its "Latest Target File" /
#     "Latest Handback File" / "Latest Handback DateTime" columns filled in
#     for both tracked docs, plus a hyperlink on the new target-file cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ---------------------------------------------------------------------
# 1. Overview sheet: status text for both tracked documents changes from
#    "Ready for handoff" to "Handed back: in sync with en-US".
# ---------------------------------------------------------------------
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws1.Range("E3").Value = "Handed back: in sync with en-US"
$ws1.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in Latest Target File / Latest Handback File /
#    Latest Handback DateTime for both rows, and link the target file.
# ---------------------------------------------------------------------
$ws2.Range("I2").Value = "73f4d9c2-45bf-46c6-be33-629ec46c93b9.md"
$ws2.Range("J2").Value = "73f4d9c2-45bf-46c6-be33-629ec46c93b9.64d8ac63ca7d46d701e3e94e8cd0b4ac6519daeb.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-04 16:52:56"

$ws2.Range("I3").Value = "eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md"
$ws2.Range("J3").Value = "eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.c43714325359d33608041ee4188111f8aa943a5e.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-09-04 16:52:56"

$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2499c35ce336dbca9590f680a136302312995a16/e2e/73f4d9c2-45bf-46c6-be33-629ec46c93b9.md", [Type]::Missing, [Type]::Missing, "73f4d9c2-45bf-46c6-be33-629ec46c93b9.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2499c35ce336dbca9590f680a136302312995a16/e2e/eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md", [Type]::Missing, [Type]::Missing, "eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md")

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape of update, different handback timestamp.
# ---------------------------------------------------------------------
$ws3.Range("I2").Value = "73f4d9c2-45bf-46c6-be33-629ec46c93b9.md"
$ws3.Range("J2").Value = "73f4d9c2-45bf-46c6-be33-629ec46c93b9.64d8ac63ca7d46d701e3e94e8cd0b4ac6519daeb.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-04 16:53:09"

$ws3.Range("I3").Value = "eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md"
$ws3.Range("J3").Value = "eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.c43714325359d33608041ee4188111f8aa943a5e.de-de.xlf"
$ws3.Range("K3").Value = "2016-09-04 16:53:09"

$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2499c35ce336dbca9590f680a136302312995a16/e2e/73f4d9c2-45bf-46c6-be33-629ec46c93b9.md", [Type]::Missing, [Type]::Missing, "73f4d9c2-45bf-46c6-be33-629ec46c93b9.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2499c35ce336dbca9590f680a136302312995a16/e2e/eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md", [Type]::Missing, [Type]::Missing, "eaeb81fe-72f6-4fc6-bd4a-12d0b6d6fcc3.md")

# ---------------------------------------------------------------------
# 4. Column widths: the longer status/strings widen a few columns
#    (best-effort; Excel re-derives these from rendered text width).
# ---------------------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 29.98
$ws1.Columns.Item(6).ColumnWidth = 29.98

$ws2.Columns.Item(3).ColumnWidth = 29.98
$ws2.Columns.Item(9).ColumnWidth = 40
$ws2.Columns.Item(10).ColumnWidth = 40

$ws3.Columns.Item(3).ColumnWidth = 29.98
$ws3.Columns.Item(9).ColumnWidth = 40
$ws3.Columns.Item(10).ColumnWidth = 40
